$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the greeting text for the 08:00-11:00 rule row (was "Good Morning").
$ws.Range("E8").Value = "GIT UPDATE"

# Match the author's cursor position left behind in the saved file.
$ws.Range("E8").Select() | Out-Null
